$d = $word.ActiveDocument

# Splits the text of a paragraph into several sibling runs (one run per
# token) instead of a single run, without Word silently re-merging the
# freshly created runs back together (same-formatting adjacent runs get
# coalesced on save). The trick: briefly bracket each split point with a
# bookmark. While the bookmark exists it forces a run boundary; once it is
# deleted the boundary is kept but the temporary bookmark is gone.
function Split-ParagraphIntoRuns($paragraph, $tokens) {
    $startPos = $paragraph.Range.Start

    # Replace the whole paragraph text with just the first token.
    $fullRange = $d.Range($startPos, $paragraph.Range.End)
    $fullRange.Text = $tokens[0]

    $bmCounter = 0
    for ($i = 1; $i -lt $tokens.Count; $i++) {
        $pEnd = $paragraph.Range.End
        $splitPos = $pEnd - 1
        $bmCounter = $bmCounter + 1
        $bmName = "SplitMark" + $startPos + "_" + $bmCounter

        $d.Bookmarks.Add($bmName, $d.Range($splitPos, $splitPos))
        $insPoint = $d.Range($splitPos, $splitPos)
        $insPoint.InsertAfter($tokens[$i])
        $d.Bookmarks($bmName).Delete()
    }
}

function Find-ParagraphByStyle($styleName, $expectedText) {
    foreach ($p in $d.Paragraphs) {
        # Paragraph.Range.Text includes the trailing paragraph mark (\r),
        # so trim it before comparing against the plain expected text.
        $paraText = $p.Range.Text.TrimEnd()
        if ($p.Style.NameLocal -eq $styleName -and $paraText -eq $expectedText) {
            return $p
        }
    }
    return $null
}

# --- Title: "Answers: Introduction to vectors" ---
$titlePara = Find-ParagraphByStyle "Title" "Answers: Introduction to vectors"
$titleTokens = @("Answers:", " ", "Introduction", " ", "to", " ", "vectors")
Split-ParagraphIntoRuns $titlePara $titleTokens

# --- Author: "Zheng Chen" ---
$authorPara = Find-ParagraphByStyle "Author" "Zheng Chen"
$authorTokens = @("Zheng", " ", "Chen")
Split-ParagraphIntoRuns $authorPara $authorTokens

# --- Abstract: "Answers to questions relating to the guide on introduction to vectors." ---
$abstractPara = Find-ParagraphByStyle "Abstract" "Answers to questions relating to the guide on introduction to vectors."
$abstractTokens = @("Answers", " ", "to", " ", "questions", " ", "relating", " ", "to", " ", "the", " ", "guide", " ", "on", " ", "introduction", " ", "to", " ", "vectors.")
Split-ParagraphIntoRuns $abstractPara $abstractTokens

Write-Host "Title:" $d.Paragraphs(1).Range.Text
Write-Host "Author:" $d.Paragraphs(2).Range.Text
Write-Host "Abstract:" $d.Paragraphs(4).Range.Text
